$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- 1) Re-write the existing first paragraph: add rPr (color FFFFFF + empty underline)
#        to both the text run and the line-break run. ---
$p1 = $d.Paragraphs(1)
$p1xml = "<w:p $wns>" +
           "<w:pPr><w:spacing w:after=`"100`"/></w:pPr>" +
           "<w:r>" +
             "<w:rPr><w:color w:val=`"FFFFFF`"/><w:u w:val=`"`"/></w:rPr>" +
             "<w:t>Eu,teste query teste query, com número de CPF teste query e RG teste query</w:t>" +
           "</w:r>" +
           "<w:r>" +
             "<w:rPr><w:color w:val=`"FFFFFF`"/><w:u w:val=`"`"/></w:rPr>" +
             "<w:br/>" +
           "</w:r>" +
         "</w:p>"
$p1.Range.InsertXML($p1xml)

# --- 2) Append three new paragraphs after it: "Modelo Número 1" + break,
#        an empty paragraph, and an empty paragraph with a bottom border. ---
$p1 = $d.Paragraphs(1)
$endRange = $d.Range($p1.Range.End, $p1.Range.End)

$newXml = "<w:p $wns>" +
            "<w:pPr><w:spacing w:after=`"100`"/></w:pPr>" +
            "<w:r><w:t>Modelo Número 1</w:t></w:r>" +
            "<w:r><w:br/></w:r>" +
          "</w:p>" +
          "<w:p $wns><w:pPr/></w:p>" +
          "<w:p $wns>" +
            "<w:pPr>" +
              "<w:pBdr><w:bottom w:val=`"single`" w:sz=`"6`" w:space=`"1`" w:color=`"auto`"/></w:pBdr>" +
            "</w:pPr>" +
          "</w:p>"

$endRange.InsertXML($newXml)

Write-Output "done"
